$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H11").Value = 12472.8
$ws.Range("I11").Value = 12472.8
$ws.Range("K11").Value = 12472.8
$ws.Range("M11").Value = -12332.8
$ws.Range("H21").Value = 0
$ws.Range("I21").Value = 0
$ws.Range("K21").Value = 0
$ws.Range("M21").ClearContents()
$ws.Range("H23").Value = 0
$ws.Range("I23").Value = 0
$ws.Range("K23").Value = 0
$ws.Range("M23").ClearContents()
$ws.Range("H112").Value = 4565.1577
$ws.Range("I112").Value = 2266.6667
$ws.Range("J112").Value = 4996.125
$ws.Range("K112").Value = 6800.000100000001
$ws.Range("L112").Value = 14988.375
$ws.Range("M112").Value = -5692.000100000001
$ws.Range("N112").Value = -17204.375
$ws.Range("H137").Value = 3704583.8
$ws.Range("I137").Value = 808.4761999999999
$ws.Range("K137").Value = 2425.4286
$ws.Range("M137").Value = 124.5714000000003

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1401
$ws.Range("I2").Value = 973.8333
$ws.Range("J2").Value = 2938.8
$ws.Range("K2").Value = 973.8333
$ws.Range("L2").Value = 2938.8
$ws.Range("M2").Value = -860.8333
$ws.Range("N2").Value = -3164.8
$ws.Range("H32").Value = 202168.6
$ws.Range("I32").Value = 244868.92
$ws.Range("K32").Value = 244868.92
$ws.Range("M32").Value = -244581.92
$ws.Range("H61").Value = 2266805.8
$ws.Range("I61").Value = 4806.1904
$ws.Range("J61").Value = 11767204
$ws.Range("K61").Value = 4806.1904
$ws.Range("L61").Value = 11767204
$ws.Range("M61").Value = -4594.1904
$ws.Range("N61").Value = -11767628
$ws.Range("H116").Value = 1401
$ws.Range("I116").Value = 973.8333
$ws.Range("J116").Value = 2938.8
$ws.Range("K116").Value = 973.8333
$ws.Range("L116").Value = 2938.8
$ws.Range("M116").Value = 1320.1667
$ws.Range("N116").Value = -7526.8
$ws.Range("H132").Value = 5640.6
$ws.Range("I132").Value = 5552
$ws.Range("K132").Value = 16656
$ws.Range("M132").Value = -14126
$ws.Range("H136").Value = 2266805.8
$ws.Range("I136").Value = 4806.1904
$ws.Range("J136").Value = 11767204
$ws.Range("K136").Value = 14418.5712
$ws.Range("L136").Value = 35301612
$ws.Range("M136").Value = -11868.5712
$ws.Range("N136").Value = -35306712

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1401
$ws.Range("I3").Value = 973.8333
$ws.Range("J3").Value = 2938.8
$ws.Range("K3").Value = 973.8333
$ws.Range("L3").Value = 2938.8
$ws.Range("M3").Value = -859.8333
$ws.Range("N3").Value = -3166.8
$ws.Range("H134").Value = 29034528
$ws.Range("I134").Value = 2389.1177
$ws.Range("J134").Value = 64287840
$ws.Range("K134").Value = 7167.353099999999
$ws.Range("L134").Value = 192863520
$ws.Range("M134").Value = -4632.353099999999
$ws.Range("N134").Value = -192868590

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2625.2964
$ws.Range("I31").Value = 2253.5173
$ws.Range("K31").Value = 2253.5173
$ws.Range("M31").Value = -1958.5173
$ws.Range("H34").Value = 2625.2964
$ws.Range("I34").Value = 2253.5173
$ws.Range("K34").Value = 2253.5173
$ws.Range("M34").Value = -2051.5173
$ws.Range("H107").Value = 1368.409
$ws.Range("I107").Value = 1233.4667
$ws.Range("K107").Value = 1233.4667
$ws.Range("M107").Value = 686.5333000000001
$ws.Range("H132").Value = 24531.777
$ws.Range("I132").Value = 32401.697
$ws.Range("K132").Value = 97205.091
$ws.Range("M132").Value = -94675.091

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H32").Value = 9991
$ws.Range("I32").Value = 9991
$ws.Range("K32").Value = 29973
$ws.Range("M32").Value = -29690
$ws.Range("H38").Value = 435.9375
$ws.Range("I38").Value = 289.66666
$ws.Range("J38").Value = 624
$ws.Range("K38").Value = 868.9999799999999
$ws.Range("L38").Value = 1872
$ws.Range("M38").Value = -521.9999799999999
$ws.Range("N38").Value = -2566
$ws.Range("H44").Value = 2950
$ws.Range("I44").Value = 2000
$ws.Range("J44").Value = 3900
$ws.Range("K44").Value = 6000
$ws.Range("L44").Value = 11700
$ws.Range("M44").Value = -5602
$ws.Range("N44").Value = -12496
$ws.Range("H70").Value = 7898
$ws.Range("I70").Value = 1849.5
$ws.Range("K70").Value = 5548.5
$ws.Range("M70").Value = -5233.5
$ws.Range("H73").Value = 7898
$ws.Range("I73").Value = 1849.5
$ws.Range("K73").Value = 5548.5
$ws.Range("M73").Value = -4456.5
$ws.Range("H81").Value = 7734.4165
$ws.Range("I81").Value = 1961.6666
$ws.Range("J81").Value = 9658.666999999999
$ws.Range("K81").Value = 5884.9998
$ws.Range("L81").Value = 28976.001
$ws.Range("M81").Value = -4761.9998
$ws.Range("N81").Value = -31222.001
$ws.Range("H84").Value = 7734.4165
$ws.Range("I84").Value = 1961.6666
$ws.Range("J84").Value = 9658.666999999999
$ws.Range("K84").Value = 17654.9994
$ws.Range("L84").Value = 86928.003
$ws.Range("M84").Value = -12038.9994
$ws.Range("N84").Value = -98160.003
$ws.Range("H113").Value = 516.7273
$ws.Range("J113").Value = 520.6667
$ws.Range("L113").Value = 1562.0001
$ws.Range("N113").Value = -5902.0001
$ws.Range("H115").Value = 1410.5555
$ws.Range("I115").Value = 385.14285
$ws.Range("K115").Value = 1155.42855
$ws.Range("M115").Value = 19.57144999999991
$ws.Range("H122").Value = 6927551.5
$ws.Range("I122").Value = 13853470
$ws.Range("K122").Value = 124681230
$ws.Range("M122").Value = -124678780
$ws.Range("H137").Value = 3427.3333
$ws.Range("I137").Value = 3414.4443
$ws.Range("J137").Value = 3466
$ws.Range("K137").Value = 10243.3329
$ws.Range("L137").Value = 10398
$ws.Range("M137").Value = -5143.332900000001
$ws.Range("N137").Value = -20598

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H10").Value = 9334.333000000001
$ws.Range("I10").Value = 5003
$ws.Range("K10").Value = 5003
$ws.Range("M10").Value = -4834
$ws.Range("H102").Value = 71431160
$ws.Range("I102").Value = 83335420
$ws.Range("K102").Value = 83335420
$ws.Range("M102").Value = -83333798
$ws.Range("H126").Value = 2871
$ws.Range("I126").Value = 2494.6667
$ws.Range("J126").Value = 4000
$ws.Range("K126").Value = 7484.000100000001
$ws.Range("L126").Value = 12000
$ws.Range("M126").Value = -5014.000100000001
$ws.Range("N126").Value = -16940
$ws.Range("H132").Value = 550773.4
$ws.Range("I132").Value = 7270.2856
$ws.Range("J132").Value = 989756.7
$ws.Range("K132").Value = 21810.8568
$ws.Range("L132").Value = 2969270.1
$ws.Range("M132").Value = -19280.8568
$ws.Range("N132").Value = -2974330.1

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H4").Value = 7009
$ws.Range("I4").Value = 7009
$ws.Range("J4").Value = 0
$ws.Range("K4").Value = 7009
$ws.Range("L4").Value = 0
$ws.Range("M4").Value = -6896
$ws.Range("N4").ClearContents()
$ws.Range("H16").Value = 1287.4849
$ws.Range("I16").Value = 1134.0344
$ws.Range("K16").Value = 1134.0344
$ws.Range("M16").Value = -964.0344
$ws.Range("H22").Value = 4742.3447
$ws.Range("I22").Value = 2291.6155
$ws.Range("K22").Value = 2291.6155
$ws.Range("M22").Value = -1996.6155
$ws.Range("H27").Value = 4742.3447
$ws.Range("I27").Value = 2291.6155
$ws.Range("K27").Value = 2291.6155
$ws.Range("M27").Value = -2184.6155
$ws.Range("H28").Value = 7009
$ws.Range("I28").Value = 7009
$ws.Range("J28").Value = 0
$ws.Range("K28").Value = 7009
$ws.Range("L28").Value = 0
$ws.Range("M28").Value = -6777
$ws.Range("N28").ClearContents()
$ws.Range("H37").Value = 7009
$ws.Range("I37").Value = 7009
$ws.Range("J37").Value = 0
$ws.Range("K37").Value = 7009
$ws.Range("L37").Value = 0
$ws.Range("M37").Value = -6902
$ws.Range("N37").ClearContents()
$ws.Range("H103").Value = 29999.5
$ws.Range("J103").Value = 29999.5
$ws.Range("L103").Value = 29999.5
$ws.Range("N103").Value = -32343.5
$ws.Range("H122").Value = 3273.3262
$ws.Range("I122").Value = 2919.1875
$ws.Range("J122").Value = 4082.7856
$ws.Range("K122").Value = 8757.5625
$ws.Range("L122").Value = 12248.3568
$ws.Range("M122").Value = -6307.5625
$ws.Range("N122").Value = -17148.3568
$ws.Range("H132").Value = 2863.5217
$ws.Range("I132").Value = 2521.2942
$ws.Range("J132").Value = 3833.1667
$ws.Range("K132").Value = 7563.882599999999
$ws.Range("L132").Value = 11499.5001
$ws.Range("M132").Value = -5033.882599999999
$ws.Range("N132").Value = -16559.5001
$ws.Range("H136").Value = 2793.9382
$ws.Range("J136").Value = 2728.1804
$ws.Range("L136").Value = 8184.541200000001
$ws.Range("N136").Value = -13284.5412

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 2640.4048
$ws.Range("J132").Value = 2982.05
$ws.Range("L132").Value = 8946.150000000001
$ws.Range("N132").Value = -14006.15
$ws.Range("H133").Value = 53873.75
$ws.Range("J133").Value = 53873.75
$ws.Range("L133").Value = 53873.75
$ws.Range("N133").Value = -63993.75
